$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 6.906113890363072
$ws.Range("C2").Value = 4.90392508821604
$ws.Range("D2").Value = 4.900273728404036
$ws.Range("F2").Value = 23.94134360223773
$ws.Range("G2").Value = 3.629648391171952
$ws.Range("I2").Value = 20.59350399651791
$ws.Range("K2").Value = 7.682836885518371
$ws.Range("M2").Value = 19.97477329853833
$ws.Range("O2").Value = 21.4559750341556
$ws.Range("B3").Value = 6.651024515813699
$ws.Range("C3").Value = 4.800359225142183
$ws.Range("D3").Value = 4.840126806367293
$ws.Range("F3").Value = 23.95451916270414
$ws.Range("G3").Value = 3.631347084308849
$ws.Range("I3").Value = 20.66926421772198
$ws.Range("K3").Value = 7.480574872430257
$ws.Range("M3").Value = 19.38065203147536
$ws.Range("O3").Value = 21.51396565792024
$ws.Range("B4").Value = 6.490047010523687
$ws.Range("C4").Value = 4.734849499013662
$ws.Range("D4").Value = 4.802192482888754
$ws.Range("F4").Value = 23.96889503291076
$ws.Range("G4").Value = 3.632444432776203
$ws.Range("I4").Value = 20.71973086359063
$ws.Range("K4").Value = 7.352066041931496
$ws.Range("M4").Value = 19.01464140048856
$ws.Range("O4").Value = 21.55412970650277
$ws.Range("B5").Value = 6.423454267693673
$ws.Range("C5").Value = 4.707692753964075
$ws.Range("D5").Value = 4.786489500542277
$ws.Range("F5").Value = 23.97633157403374
$ws.Range("G5").Value = 3.632905323231519
$ws.Range("I5").Value = 20.74128777818588
$ws.Range("K5").Value = 7.298656395021336
$ws.Range("M5").Value = 18.86543315416207
$ws.Range("O5").Value = 21.57163992270001
$ws.Range("B6").Value = 6.412339774128373
$ws.Range("C6").Value = 4.703156214394017
$ws.Range("D6").Value = 4.783867506743624
$ws.Range("F6").Value = 23.97766164113933
$ws.Range("G6").Value = 3.632982683288928
$ws.Range("I6").Value = 20.74492710958719
$ws.Range("K6").Value = 7.289726254050535
$ws.Range("M6").Value = 18.84066140130647
$ws.Range("O6").Value = 21.57461642715657
$ws.Range("B7").Value = 6.489152797170668
$ws.Range("C7").Value = 4.734485090580154
$ws.Range("D7").Value = 4.801981685999056
$ws.Range("F7").Value = 23.96898893831
$ws.Range("G7").Value = 3.632450592924238
$ws.Range("I7").Value = 20.7200175762055
$ws.Range("K7").Value = 7.351349894976162
$ws.Range("M7").Value = 19.01262901319661
$ws.Range("O7").Value = 21.55436123140719
$ws.Range("B8").Value = 6.819119710679633
$ws.Range("C8").Value = 4.868626714087617
$ws.Range("D8").Value = 4.8797481453649
$ws.Range("F8").Value = 23.94458067614854
$ws.Range("G8").Value = 3.630222849709071
$ws.Range("I8").Value = 20.61880530723138
$ws.Range("K8").Value = 7.614019175311683
$ws.Range("M8").Value = 19.77031937119033
$ws.Range("O8").Value = 21.47502251551109
$ws.Range("B9").Value = 7.427733355992561
$ws.Range("C9").Value = 5.115553971211544
$ws.Range("D9").Value = 5.023944755905126
$ws.Range("F9").Value = 23.94666422175728
$ws.Range("G9").Value = 3.626283305734864
$ws.Range("I9").Value = 20.45174684640185
$ws.Range("K9").Value = 8.093053631636492
$ws.Range("M9").Value = 21.23615979068745
$ws.Range("O9").Value = 21.35573091182922
$ws.Range("B10").Value = 7.846878341013432
$ws.Range("C10").Value = 5.286041696436227
$ws.Range("D10").Value = 5.124391205843716
$ws.Range("F10").Value = 23.97868821206652
$ws.Range("G10").Value = 3.623647478601755
$ws.Range("I10").Value = 20.34827251236088
$ws.Range("K10").Value = 8.420868240447209
$ws.Range("M10").Value = 22.28817525209275
$ws.Range("O10").Value = 21.2903806224609
$ws.Range("B11").Value = 8.030695325010484
$ws.Range("C11").Value = 5.361020914926681
$ws.Range("D11").Value = 5.168807358241391
$ws.Range("F11").Value = 23.99986773328408
$ws.Range("G11").Value = 3.622503872656329
$ws.Range("I11").Value = 20.30540591654382
$ws.Range("K11").Value = 8.564371450850055
$ws.Range("M11").Value = 22.75898319465318
$ws.Range("O11").Value = 21.26552481857989
$ws.Range("B12").Value = 8.099262113516405
$ws.Range("C12").Value = 5.389027490351479
$ws.Range("D12").Value = 5.18543598796014
$ws.Range("F12").Value = 24.0088360179324
$ws.Range("G12").Value = 3.622078742653709
$ws.Range("I12").Value = 20.28978004655194
$ws.Range("K12").Value = 8.617875710998794
$ws.Range("M12").Value = 22.93597616957267
$ws.Range("O12").Value = 21.25681561874956
$ws.Range("B13").Value = 8.084542115063096
$ws.Range("C13").Value = 5.383013200007989
$ws.Range("D13").Value = 5.181863311021915
$ws.Range("F13").Value = 24.00686242562033
$ws.Range("G13").Value = 3.622169950074373
$ws.Range("I13").Value = 20.29311833569549
$ws.Range("K13").Value = 8.60639028006791
$ws.Range("M13").Value = 22.8979176384697
$ws.Range("O13").Value = 21.25865998987605
$ws.Range("B14").Value = 8.036357503167654
$ws.Range("C14").Value = 5.363332864111706
$ws.Range("D14").Value = 5.170179269244326
$ws.Range("F14").Value = 24.00058657789826
$ws.Range("G14").Value = 3.622468738303317
$ws.Range("I14").Value = 20.30410819351872
$ws.Range("K14").Value = 8.568790214577819
$ws.Range("M14").Value = 22.77357134967487
$ws.Range("O14").Value = 21.26479420153015
$ws.Range("B15").Value = 8.006705994797294
$ws.Range("C15").Value = 5.351227314514706
$ws.Range("D15").Value = 5.162997402236521
$ws.Range("F15").Value = 23.99686580839913
$ws.Range("G15").Value = 3.622652785969832
$ws.Range("I15").Value = 20.31091888431558
$ws.Range("K15").Value = 8.545649239329714
$ws.Range("M15").Value = 22.69723249859559
$ws.Range("O15").Value = 21.26864322587476
$ws.Range("B16").Value = 7.834722667394659
$ws.Range("C16").Value = 5.281088391914172
$ws.Range("D16").Value = 5.121462185908329
$ws.Range("F16").Value = 23.97743691702175
$ws.Range("G16").Value = 3.623723327854672
$ws.Range("I16").Value = 20.35115873596635
$ws.Range("K16").Value = 8.411374310784545
$ws.Range("M16").Value = 22.2572352731833
$ws.Range("O16").Value = 21.29210327148666
$ws.Range("B17").Value = 7.727419383591363
$ws.Range("C17").Value = 5.2373887959057
$ws.Range("D17").Value = 5.095649339707123
$ws.Range("F17").Value = 23.96720972986696
$ws.Range("G17").Value = 3.624394239940824
$ws.Range("I17").Value = 20.37692299012333
$ws.Range("K17").Value = 8.32753969868102
$ws.Range("M17").Value = 21.98519322277521
$ws.Range("O17").Value = 21.30774502978325
$ws.Range("B18").Value = 7.665059024187662
$ws.Range("C18").Value = 5.21201230260514
$ws.Range("D18").Value = 5.080682640740085
$ws.Range("F18").Value = 23.96194986665316
$ws.Range("G18").Value = 3.624785352209274
$ws.Range("I18").Value = 20.39213751482774
$ws.Range("K18").Value = 8.278792681372613
$ws.Range("M18").Value = 21.82800092531548
$ws.Range("O18").Value = 21.31720022798915
$ws.Range("B19").Value = 7.643836354196348
$ws.Range("C19").Value = 5.203379261867945
$ws.Range("D19").Value = 5.075594813788021
$ws.Range("F19").Value = 23.96027595527953
$ws.Range("G19").Value = 3.624918674228697
$ws.Range("I19").Value = 20.3973567704297
$ws.Range("K19").Value = 8.26219809967607
$ws.Range("M19").Value = 21.77466014298068
$ws.Range("O19").Value = 21.32048025893306
$ws.Range("B20").Value = 7.738908904370439
$ws.Range("C20").Value = 5.242065820855109
$ws.Range("D20").Value = 5.098409615912031
$ws.Range("F20").Value = 23.96823401671524
$ws.Range("G20").Value = 3.624322280131971
$ws.Range("I20").Value = 20.37413938156313
$ws.Range("K20").Value = 8.33651885637652
$ws.Range("M20").Value = 22.01422850394869
$ws.Range("O20").Value = 21.30603246893731
$ws.Range("B21").Value = 8.050539139415593
$ws.Range("C21").Value = 5.369124064084395
$ws.Range("D21").Value = 5.173616389039815
$ws.Range("F21").Value = 24.00240424282879
$ws.Range("G21").Value = 3.622380762125106
$ws.Range("I21").Value = 20.30086372061129
$ws.Range("K21").Value = 8.579857205378664
$ws.Range("M21").Value = 22.81013124521196
$ws.Range("O21").Value = 21.2629733324273
$ws.Range("B22").Value = 8.248120861420391
$ws.Range("C22").Value = 5.449905063166895
$ws.Range("D22").Value = 5.221653132013253
$ws.Range("F22").Value = 24.03026052279678
$ws.Range("G22").Value = 3.621158062303281
$ws.Range("I22").Value = 20.25651162802919
$ws.Range("K22").Value = 8.734001100697517
$ws.Range("M22").Value = 23.32270058461961
$ws.Range("O22").Value = 21.23893124079337
$ws.Range("B23").Value = 8.143241062054695
$ws.Range("C23").Value = 5.407002354537456
$ws.Range("D23").Value = 5.196119332211639
$ws.Range("F23").Value = 24.01488879693613
$ws.Range("G23").Value = 3.621806427969219
$ws.Range("I23").Value = 20.27985874400523
$ws.Range("K23").Value = 8.652188021020834
$ws.Range("M23").Value = 23.04988174876965
$ws.Range("O23").Value = 21.25138706988315
$ws.Range("B24").Value = 7.733716575517807
$ws.Range("C24").Value = 5.239952124733923
$ws.Range("D24").Value = 5.097162088672271
$ws.Range("F24").Value = 23.9677690056422
$ws.Range("G24").Value = 3.624354796351112
$ws.Range("I24").Value = 20.37539659804779
$ws.Range("K24").Value = 8.332461088157011
$ws.Range("M24").Value = 22.00110410575252
$ws.Range("O24").Value = 21.30680527717003
$ws.Range("B25").Value = 7.267711134143698
$ws.Range("C25").Value = 5.050598443924152
$ws.Range("D25").Value = 4.985870384122716
$ws.Range("F25").Value = 23.940745006513
$ws.Range("G25").Value = 3.627303435118946
$ws.Range("I25").Value = 20.49356566438687
$ws.Range("K25").Value = 7.967557510412552
$ws.Range("M25").Value = 20.8431302997314
$ws.Range("O25").Value = 21.38409910152391
